# Update "想去人数" (want-to-go count) figures with freshly scraped values.
$wb = $excel.ActiveWorkbook

# Sheet "展览"
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F4").Value = 142
$ws1.Range("F8").Value = 1112
$ws1.Range("F12").Value = 648
$ws1.Range("F26").Value = 2435
$ws1.Range("F29").Value = 14
$ws1.Range("F31").Value = 267

# Sheet "演出"
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F8").Value = 228
$ws2.Range("F27").Value = 3822

# Sheet "全部类型"
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F10").Value = 142
$ws4.Range("F15").Value = 1112
$ws4.Range("F18").Value = 648
$ws4.Range("F31").Value = 2435
$ws4.Range("F36").Value = 14
$ws4.Range("F39").Value = 267
